# Apply F-column ("想去人数" / want-to-go count) updates across the
# four sheets (展览 / 演出 / 本地生活 / 全部类型) to match the regenerated
# data pull (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# Sheet 1: 展览
$ws1.Range("F2").Value = 7677
$ws1.Range("F9").Value = 197
$ws1.Range("F24").Value = 2673
$ws1.Range("F25").Value = 129
$ws1.Range("F26").Value = 118
$ws1.Range("F27").Value = 3142
$ws1.Range("F28").Value = 2443
$ws1.Range("F29").Value = 79
$ws1.Range("F34").Value = 144
$ws1.Range("F35").Value = 50
$ws1.Range("F36").Value = 41
$ws1.Range("F38").Value = 4593
$ws1.Range("F39").Value = 567
$ws1.Range("F43").Value = 893
$ws1.Range("F44").Value = 292
$ws1.Range("F45").Value = 16

# Sheet 2: 演出
$ws2.Range("F7").Value = 70

# Sheet 3: 本地生活
$ws3.Range("F3").Value = 11

# Sheet 4: 全部类型
$ws4.Range("F4").Value = 7677
$ws4.Range("F11").Value = 197
$ws4.Range("F13").Value = 70
$ws4.Range("F25").Value = 2673
$ws4.Range("F26").Value = 129
$ws4.Range("F27").Value = 118
$ws4.Range("F28").Value = 2443
$ws4.Range("F29").Value = 79
$ws4.Range("F32").Value = 144
$ws4.Range("F33").Value = 50
$ws4.Range("F34").Value = 41
$ws4.Range("F38").Value = 4593
$ws4.Range("F40").Value = 567
$ws4.Range("F44").Value = 893
$ws4.Range("F45").Value = 292
